$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a flat data table (rows 2..435) of weekly price observations,
# two rows per date (quality "Primera" then "Segunda"). A new week's
# observations (date serial 44511) are inserted at the very top of the
# "Apio" block (rows 349-350), pushing the rest of the table down by two
# rows (old row 349 -> new row 351, ..., old row 435 -> new row 437).

$ws.Range("A349:A350").EntireRow.Insert()

# New row 349 - Calidad "Primera"
$ws.Cells.Item(349, 1).Value = 6
$ws.Cells.Item(349, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(349, 3).Value = "Metropolitana"
$ws.Cells.Item(349, 4).Value = 44511
$ws.Cells.Item(349, 5).Value = 13
$ws.Cells.Item(349, 6).Value = 100112017
$ws.Cells.Item(349, 7).Value = "Apio"
$ws.Cells.Item(349, 8).Value = "Americana (o)"
$ws.Cells.Item(349, 9).Value = "Primera"
$ws.Cells.Item(349, 10).Value = 2200
$ws.Cells.Item(349, 11).Value = 6000
$ws.Cells.Item(349, 12).Value = 7000
$ws.Cells.Item(349, 13).Value = 6568
$ws.Cells.Item(349, 14).Value = "`$/docena de matas"
$ws.Cells.Item(349, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(349, 16).Value = 1095
$ws.Cells.Item(349, 17).Value = 6
$ws.Cells.Item(349, 18).Value = "Hortaliza"

# New row 350 - Calidad "Segunda"
$ws.Cells.Item(350, 1).Value = 6
$ws.Cells.Item(350, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(350, 3).Value = "Metropolitana"
$ws.Cells.Item(350, 4).Value = 44511
$ws.Cells.Item(350, 5).Value = 13
$ws.Cells.Item(350, 6).Value = 100112017
$ws.Cells.Item(350, 7).Value = "Apio"
$ws.Cells.Item(350, 8).Value = "Americana (o)"
$ws.Cells.Item(350, 9).Value = "Segunda"
$ws.Cells.Item(350, 10).Value = 750
$ws.Cells.Item(350, 11).Value = 5000
$ws.Cells.Item(350, 12).Value = 5000
$ws.Cells.Item(350, 13).Value = 5000
$ws.Cells.Item(350, 14).Value = "`$/docena de matas"
$ws.Cells.Item(350, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(350, 16).Value = 833
$ws.Cells.Item(350, 17).Value = 6
$ws.Cells.Item(350, 18).Value = "Hortaliza"
